# Auto-generated edit script: updates Leve profit calculation columns (H-N)
# across multiple crafting-job sheets, per scheduled price-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 3196.9443
$ws.Range("I43").Value = 930
$ws.Range("J43").Value = 6030.625
$ws.Range("K43").Value = 930
$ws.Range("L43").Value = 6030.625
$ws.Range("M43").Value = -861
$ws.Range("N43").Value = -6168.625

$ws.Range("H112").Value = 1281.3115
$ws.Range("J112").Value = 1281.3115
$ws.Range("L112").Value = 3843.9345
$ws.Range("N112").Value = -6059.9345

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 1451.6364
$ws.Range("I61").Value = 1548.8889
$ws.Range("K61").Value = 1548.8889
$ws.Range("M61").Value = -1336.8889

$ws.Range("H74").Value = 1505.1082
$ws.Range("I74").Value = 1036.9395
$ws.Range("J74").Value = 5367.5
$ws.Range("K74").Value = 1036.9395
$ws.Range("L74").Value = 5367.5
$ws.Range("M74").Value = -162.9395
$ws.Range("N74").Value = -7115.5

$ws.Range("H77").Value = 1505.1082
$ws.Range("I77").Value = 1036.9395
$ws.Range("J77").Value = 5367.5
$ws.Range("K77").Value = 5184.6975
$ws.Range("L77").Value = 26837.5
$ws.Range("M77").Value = -816.6975000000002
$ws.Range("N77").Value = -35573.5

$ws.Range("H122").Value = 5657.048
$ws.Range("I122").Value = 4906.4614
$ws.Range("J122").Value = 6876.75
$ws.Range("K122").Value = 14719.3842
$ws.Range("L122").Value = 20630.25
$ws.Range("M122").Value = -12269.3842
$ws.Range("N122").Value = -25530.25

$ws.Range("H130").Value = 50429
$ws.Range("J130").Value = 50429
$ws.Range("L130").Value = 50429
$ws.Range("N130").Value = -60469

$ws.Range("H136").Value = 1451.6364
$ws.Range("I136").Value = 1548.8889
$ws.Range("K136").Value = 4646.6667
$ws.Range("M136").Value = -2096.6667

$ws.Range("H137").Value = 45774
$ws.Range("J137").Value = 45774
$ws.Range("L137").Value = 45774
$ws.Range("N137").Value = -55974

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H59").Value = 38890
$ws.Range("J59").Value = 38890
$ws.Range("L59").Value = 38890
$ws.Range("N59").Value = -40584

$ws.Range("H126").Value = 41891.11
$ws.Range("J126").Value = 41891.11
$ws.Range("L126").Value = 41891.11
$ws.Range("N126").Value = -51771.11

$ws.Range("H134").Value = 2238.5217
$ws.Range("I134").Value = 1375.0555
$ws.Range("J134").Value = 5347
$ws.Range("K134").Value = 4125.166499999999
$ws.Range("L134").Value = 16041
$ws.Range("M134").Value = -1590.166499999999
$ws.Range("N134").Value = -21111

$ws.Range("H137").Value = 46075.8
$ws.Range("J137").Value = 46075.8
$ws.Range("L137").Value = 46075.8
$ws.Range("N137").Value = -56275.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 7054.3696
$ws.Range("I31").Value = 3498.4614
$ws.Range("J31").Value = 8455.182000000001
$ws.Range("K31").Value = 3498.4614
$ws.Range("L31").Value = 8455.182000000001
$ws.Range("M31").Value = -3203.4614
$ws.Range("N31").Value = -9045.182000000001

$ws.Range("H34").Value = 7054.3696
$ws.Range("I34").Value = 3498.4614
$ws.Range("J34").Value = 8455.182000000001
$ws.Range("K34").Value = 3498.4614
$ws.Range("L34").Value = 8455.182000000001
$ws.Range("M34").Value = -3296.4614
$ws.Range("N34").Value = -8859.182000000001

$ws.Range("H58").Value = 2065.111
$ws.Range("I58").Value = 1305.45
$ws.Range("J58").Value = 4235.5713
$ws.Range("K58").Value = 1305.45
$ws.Range("L58").Value = 4235.5713
$ws.Range("M58").Value = -1102.45
$ws.Range("N58").Value = -4641.5713

$ws.Range("H94").Value = 1640.4445
$ws.Range("J94").Value = 1825.9231
$ws.Range("L94").Value = 1825.9231
$ws.Range("N94").Value = -2727.9231

$ws.Range("H118").Value = 28990
$ws.Range("J118").Value = 28990
$ws.Range("L118").Value = 28990
$ws.Range("N118").Value = -32304

$ws.Range("H136").Value = 2065.111
$ws.Range("I136").Value = 1305.45
$ws.Range("J136").Value = 4235.5713
$ws.Range("K136").Value = 3916.35
$ws.Range("L136").Value = 12706.7139
$ws.Range("M136").Value = -1366.35
$ws.Range("N136").Value = -17806.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 844.04
$ws.Range("I131").Value = 320
$ws.Range("J131").Value = 871.62103
$ws.Range("K131").Value = 960
$ws.Range("L131").Value = 2614.86309
$ws.Range("M131").Value = 4080
$ws.Range("N131").Value = -12694.86309

$ws.Range("H140").Value = 3163.0625
$ws.Range("I140").Value = 4045.4443
$ws.Range("J140").Value = 2028.5714
$ws.Range("K140").Value = 12136.3329
$ws.Range("L140").Value = 6085.7142
$ws.Range("M140").Value = -6956.332900000001
$ws.Range("N140").Value = -16445.7142

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H46").Value = 26892.77
$ws.Range("J46").Value = 27029.666
$ws.Range("L46").Value = 27029.666
$ws.Range("N46").Value = -27341.666

$ws.Range("H102").Value = 3076.7368
$ws.Range("I102").Value = 1904.4286
$ws.Range("J102").Value = 6359.2
$ws.Range("K102").Value = 1904.4286
$ws.Range("L102").Value = 6359.2
$ws.Range("M102").Value = -282.4286
$ws.Range("N102").Value = -9603.200000000001

$ws.Range("H126").Value = 3412.3635
$ws.Range("I126").Value = 2828.169
$ws.Range("J126").Value = 4893.7144
$ws.Range("K126").Value = 8484.507
$ws.Range("L126").Value = 14681.1432
$ws.Range("M126").Value = -6014.507
$ws.Range("N126").Value = -19621.1432

$ws.Range("H137").Value = 43746
$ws.Range("J137").Value = 43746
$ws.Range("L137").Value = 43746
$ws.Range("N137").Value = -53946

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 1414.1724
$ws.Range("I46").Value = 835.2353000000001
$ws.Range("J46").Value = 2234.3333
$ws.Range("K46").Value = 835.2353000000001
$ws.Range("L46").Value = 2234.3333
$ws.Range("M46").Value = -647.2353000000001
$ws.Range("N46").Value = -2610.3333

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 1800
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 1800
$ws.Range("M61").Value = -1998
$ws.Range("N61").Value = -2204

$ws.Range("H93").Value = 5850049
$ws.Range("I93").Value = 9261119
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 9261119
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -9259871
$ws.Range("N93").Value = -4996

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -30
$ws.Range("N113").Value = -6140

$ws.Range("H140").Value = 59696.625
$ws.Range("J140").Value = 59696.625
$ws.Range("L140").Value = 59696.625
$ws.Range("N140").Value = -70056.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H123").Value = 38939.375
$ws.Range("J123").Value = 38939.375
$ws.Range("L123").Value = 38939.375
$ws.Range("N123").Value = -48739.375
